$wb = $excel.ActiveWorkbook

# Duplicate the most recent ranking sheet (2025-08-04) to create the new
# week's sheet (2025-08-11), placed right after it, then rename it.
$srcSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $srcSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2025-08-11"

$data = @(
    ,@(1, 'ワンパンマン', '原作/ＯＮＥ 作画/村田雄介', '208撃目')
    ,@(2, '魔王の俺が奴隷エルフを嫁にしたんだが、どう愛でればいい？', '原作／手島史詞 キャラクター原案／COMTA 漫画／板垣ハコ', '第72話')
    ,@(3, '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。', 'マツモトケンゴ', '【第７巻発売記念】描き下ろしイラスト公開！')
    ,@(4, '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―', '光永康則', '第６８話『施錠停止』⓵')
    ,@(5, '転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～', 'zunta(作画) はらわたさいぞう(原作)', '第31話：完全なる死角①')
    ,@(6, '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが', '絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)', '第4話 前編')
    ,@(7, '悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる', 'こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)', '第17話-1：「違法奴隷商討伐」')
    ,@(8, '王子様の友達', 'すけろく(著者)', '第29話')
    ,@(9, 'いとこのこ', 'いぬちく(著者)', '連載２周年記念マンガ')
    ,@(10, '元・世界１位のサブキャラ育成日記 ～廃プレイヤー、異世界を攻略中！～', '沢村治太郎(原作) 前田理想(漫画) まろ(キャラクター原案)', '第77話その1')
    ,@(11, '実は俺、最強でした？', '原作：澄守 彩 漫画：高橋 愛', '第121話　王都の危機を守れ!!・後編')
    ,@(12, 'このヒーラー、めんどくさい', '丹念に発酵(著者)', '第89話：盗賊再び')
    ,@(13, 'ダークサモナーとデキている', '車王(著者)', '第73話')
    ,@(14, 'クセ強彼女は床にいざなう', '須河篤志(著者)', '第14話後半')
    ,@(15, '勇者に全部奪われた俺は勇者の母親とパーティを組みました！', '久遠まこと(著者) 石のやっさん(原作)', '第29話')
    ,@(16, '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～', '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり', '第５１話　英雄を倒す器用貧乏（２）')
    ,@(17, '異世界魔王と召喚少女の奴隷魔術', '原作：むらさきゆきや 漫画：福田直叶 キャラクター原案：鶴崎貴大', '第127話　戦争を終わらせてみるⅢ（前編）')
    ,@(18, '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する', '作画：マエD 原作：新人', '第5話(3)')
    ,@(19, 'まんきつしたい常連さん', 'しんみりん(著者)', '第46話後編')
    ,@(20, '望まぬ不死の冒険者', '中曽根ハイジ（漫画） 丘野 優（原作） じゃいあん（キャラクター原案）', '第59話　ヴィステルヤ（前編）')
    ,@(21, '怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった', '菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)', '第11話')
    ,@(22, '最強の少年聖騎士、転生者を狩る', '作画：御塩 原作：宇奈木ユラ', '第6話(3)')
    ,@(23, 'みつばものがたり 呪いの少女と死の輪舞《ロンド》', '堤りん(漫画) 七沢またり(原作) EURA(キャラクター原案)', '第11話：勝利の美酒')
    ,@(24, '最強勇者パーティーは愛が知りたい', '山田肌襦袢', '第28話「最後はこぶしがあればいい」')
    ,@(25, '淫獄団地', '搾精研究所(原作) 丈山雄為(漫画)', '第49話（後編）')
    ,@(26, '剥かせて！竜ケ崎さん', '一智和智', '大学生編 第13話')
    ,@(27, 'リビルドワールド', '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)', '第72話①')
    ,@(28, 'ぽんドロイド！ はまさん', 'はれやまはれぞう(著者)', '第6話')
    ,@(29, 'アザミヤコを好きになる', 'ユニティコング(原作) ツノニガウ(作画)', '第9話前編')
    ,@(30, 'バキ外伝 烈海王は異世界転生しても一向にかまわんッッ', '板垣恵介 猪原賽 陸井栄史', '第77話　応援(エール)')
    ,@(31, '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜', '戸賀 環 坂木持丸 riritto', '第50話②　祝われた家を探索してみた')
    ,@(32, 'よくわからないけれど異世界に転生していたようです', '内々けやき あし カオミン', '第137話 よくわからないけれど脱出するみたいです（１）')
    ,@(33, '異世界メイドの三ツ星グルメ ～現代ごはん作ったら王宮で大バズリしました～', 'モリタ Ｕ４ nima', '第13話（１）　ドーナツの騎士様とお土産スイーツ（１）')
    ,@(34, '聖者無双', '漫画：秋風緋色 原作：ブロッコリーライオン キャラクター原案：sime', '第91話　邂逅（前半）')
    ,@(35, '濁る瞳で何を願う ハイセルク戦記', 'トルトネン 創-taro 斎藤八呑', '第32話 轍')
    ,@(36, 'ライドンキング', '馬場康誌', '第81話 大統領と失われた神器（後編）')
    ,@(37, '美人女上司滝沢さん', 'やんBARU(著者)', '第202.5話')
    ,@(38, '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~', '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)', '第81話その2')
    ,@(39, '落ちこぼれだった兄が実は最強 ～史上最強の勇者は転生し、学園で無自覚に無双する～', '村上よしゆき 茨木野 あるてら', '第４１話　勇者、人魚王国を救い、歓迎される。あと、六邪神将が、全員来る（２）')
    ,@(40, '解雇された暗黒兵士(30代)のスローなセカンドライフ', '岡沢六十四 るれくちぇ sage・ジョー', '第71話(後編) ダリエルVS.滾り')
    ,@(41, 'バキ外伝　ガイアとシコルスキー　～ときどきノムラ 二人だけど三人暮らし～', '板垣恵介 林たかあき', '第51話 春の大敵')
    ,@(42, '婚約者に裏切られた錬金術師は、独立して『ざまぁ』します　コミック版', '漫画/すたひろ 原作/Y.A', 'chapter67【35話②】')
    ,@(43, 'アイドル辞めるけど結婚してくれますか!?', '三吉汐美(著者)', '休載イラスト')
    ,@(44, '魔のものたちは企てる', '加藤拓弐(原作) ガしガし(作画)', 'コミックス告知')
    ,@(45, '俺以外誰も採取できない素材なのに「素材採取率が低い」とパワハラする幼馴染錬金術師と絶縁した専属魔導士、辺境の町でスローライフを送りたい。', '狐御前(原作) 西岡知三(作画) ＮＯＣＯ(キャラクター原案)', '第24話-1')
    ,@(46, 'ハズレ枠の【状態異常スキル】で最強になった俺がすべてを蹂躙するまで', '鵜吉しょう（作画） 内々けやき（構成） 篠崎 芳（原作） KWKM（キャラクター原案）', '第56話　十河綾香')
    ,@(47, '回復術士のやり直し', '月夜涙(原作) 羽賀ソウケン(漫画) しおこんぶ(キャラクター原案)', '第72話-2')
    ,@(48, 'じつは義妹でした。～最近できた義理の弟の距離感がやたら近いわけ～', '堺しょうきち(著者) 白井ムク(原作) 千種みのり(キャラクター原案)', '第36話-2')
    ,@(49, '塔の管理をしてみよう', '盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)', '第91話後編')
    ,@(50, '陰キャの俺が席替えでS級美少女に囲まれたら秘密の関係が始まった。', '星野 星野(原作) バラマツヒトミ(漫画) 黒兎 ゆう(キャラクター原案)', '第4話')
)

foreach ($row in $data) {
    $r = [int]$row[0] + 1
    $newSheet.Cells.Item($r, 1).Value = [int]$row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
}
